# Generate Report for handoff
# - New source markdown GUID: 3d7a97e8-37f9-49f3-852b-8681dab3bbc0 -> 11df9c05-1057-403e-88be-89967970b575
# - New handoff hash: cc935fac8bcabfef3627a6a6cd1ab5208981a7e0 -> 10ab152c1d763908b9019c13d3244b689c84e377
# - Updated handoff timestamps
# - The old "7ac6053b...md" (Handoff transform failed) row is dropped: its row's
#   values/hyperlink-display get overwritten with the former last row's content
#   (".localization-config" / "Not to be localized"), and the now-duplicate last
#   row is removed.

$wb = $excel.ActiveWorkbook

$oldGuid = "3d7a97e8-37f9-49f3-852b-8681dab3bbc0"
$newGuid = "11df9c05-1057-403e-88be-89967970b575"
$oldHash = "cc935fac8bcabfef3627a6a6cd1ab5208981a7e0"
$newHash = "10ab152c1d763908b9019c13d3244b689c84e377"

$newMdName = "$newGuid.md"
$newZhXlf  = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf  = "$newGuid.$newHash.de-de.xlf"

# ----------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$url1A2 = "https://github.com/OpenLocalizationTest/oltest/blob/69787ff6a570a42b6c71eb0caf8021d95c94372c/e2e/$oldGuid.md"
$url1A3 = "https://github.com/OpenLocalizationTest/oltest/blob/69787ff6a570a42b6c71eb0caf8021d95c94372c/e2e/7ac6053b-ae52-4416-badd-f3e7d8b60527.md"

# Row 3 becomes what row 4 used to hold; row 4 then disappears.
$ws1.Range("B3").Value = "Not to be localized"
$ws1.Range("C3").Value = "Not to be localized"
$ws1.Range("A4:C4").ClearContents()

$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $url1A2, "", "", $newMdName) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $url1A3, "", "", ".localization-config") | Out-Null

$ws1.Rows.Item(4).Delete()

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$url2A2 = "https://github.com/OpenLocalizationTest/oltest/blob/69787ff6a570a42b6c71eb0caf8021d95c94372c/e2e/$oldGuid.md"
$url2C2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bdd18fa54e1e3efb0de0626b1c18153ff92a798a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$oldGuid.$oldHash.zh-cn.xlf"
$url2A3 = "https://github.com/OpenLocalizationTest/oltest/blob/69787ff6a570a42b6c71eb0caf8021d95c94372c/e2e/7ac6053b-ae52-4416-badd-f3e7d8b60527.md"

$ws2.Range("D2").Value = "2016-01-25 11:12:43"

$ws2.Range("B3").Value = "Not to be localized"
$ws2.Range("A4:I4").ClearContents()

$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $url2A2, "", "", $newMdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $url2C2, "", "", $newZhXlf) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $url2A3, "", "", ".localization-config") | Out-Null

$ws2.Rows.Item(4).Delete()

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$url3A2 = "https://github.com/OpenLocalizationTest/oltest/blob/69787ff6a570a42b6c71eb0caf8021d95c94372c/e2e/$oldGuid.md"
$url3C2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a7c9f142bd82fcd539c46b9689b44777466e529/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$oldGuid.$oldHash.de-de.xlf"
$url3A3 = "https://github.com/OpenLocalizationTest/oltest/blob/69787ff6a570a42b6c71eb0caf8021d95c94372c/e2e/7ac6053b-ae52-4416-badd-f3e7d8b60527.md"

$ws3.Range("D2").Value = "2016-01-25 11:12:52"

$ws3.Range("B3").Value = "Not to be localized"
$ws3.Range("A4:I4").ClearContents()

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $url3A2, "", "", $newMdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $url3C2, "", "", $newDeXlf) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $url3A3, "", "", ".localization-config") | Out-Null

$ws3.Rows.Item(4).Delete()

Write-Output "Report regenerated for handoff"
